# Latest generated outputs 2025-09-18
# Insert a new row for "Trees on adjacent land" beneath "Trees on site"
# (row 151) within the "Trees and hedges information" section, shifting
# every subsequent row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 152 - everything from 152 downward shifts to 153+
$ws.Rows("152").Insert()

# Row 151 ("Trees and hedges information" / A151,B151 stay as-is) gains the
# first field: "Trees on site"
$ws.Range("C151").Value = "Trees on site"
$ws.Range("G151").Value = "Whether trees or hedges are present on the proposed development site"
$ws.Range("H151").Value = "boolean"
$ws.Range("I151").Value = "MUST"

# New row 152: second field of the same group - "Trees on adjacent land"
# (D152/E152/F152 are left blank - same as every other unused column in
# this sub-row layout)
$ws.Range("C152").Value = "Trees on adjacent land"
$ws.Range("G152").Value = "Whether trees or hedges on land adjacent to the proposed development site could influence the development or might be important as part of the local landscape character"
$ws.Range("H152").Value = "boolean"
$ws.Range("I152").Value = "MUST"

# Extend the A151/B151 merges down to cover the new row 152
$ws.Range("A151:A152").Merge()
$ws.Range("B151:B152").Merge()
